# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds per-game strikeout values. The save_data
# generation routine was rerun (different std/mean derivation) and the
# resulting s_vals were rewritten back into column G for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => new K (strikeout) value
$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 1
    6  = 2
    7  = 2
    8  = 2
    9  = 0
    10 = 3
    11 = 2
    12 = 3
    13 = 1
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 3
    21 = 1
    22 = 1
    23 = 3
    24 = 2
    25 = 1
    26 = 2
    27 = 3
    28 = 0
    29 = 1
    30 = 2
    31 = 3
    32 = 0
    33 = 2
    34 = 2
    35 = 1
    36 = 3
    37 = 1
    38 = 3
    40 = 2
    41 = 0
    42 = 1
    44 = 0
    45 = 0
    46 = 0
    47 = 2
    48 = 1
    49 = 2
    50 = 0
    51 = 1
    52 = 1
    53 = 0
    54 = 1
    55 = 2
    56 = 2
    57 = 2
    58 = 0
    59 = 1
    60 = 0
    61 = 0
    62 = 2
    63 = 1
    64 = 1
    65 = 2
    66 = 2
    67 = 1
    68 = 2
    69 = 2
    70 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
